$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.572298049926758
$ws.Range("B1").Value = 1.637515068054199
$ws.Range("C1").Value = 1.837121605873108
$ws.Range("D1").Value = 2.82489275932312
$ws.Range("E1").Value = 3.287432432174683
